$wb = $excel.ActiveWorkbook

# --- Sheet3 content updates -----------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# New values: write U3 first so "0.66231 with all features" becomes the
# first newly-added shared string, then J3 so its text becomes the next one.
$ws3.Range("U3").Value = "0.66231 with all features"
$ws3.Range("J3").Value = "0.64388 with 15 features, 0.65082 with all features"


# --- Add Sheet4 at the end ------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"

$ws4.Range("A1").Value = "2017-08-04 09:25:34.848946 Training Random Forest classifier, [5, 2, 5, 120]"
$ws4.Range("A2").Value = "2017-08-04 09:26:45.107463 Model Training Complete"
$ws4.Range("A3").Value = "2017-08-04 09:27:10.012725 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A4").Value = "0,Train,0.946772409143,[[280717, 14552], [16865, 278105]],0.987168984168"
$ws4.Range("A5").Value = "2017-08-04 09:27:18.949925 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A6").Value = "0,Test,0.946123702013,[[93353, 4871], [5729, 92794]],0.986979032061"
$ws4.Range("A8").Value = "2017-08-04 09:28:32.729770 Model Training Complete"
$ws4.Range("A9").Value = "2017-08-04 09:28:57.777968 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A10").Value = "1,Train,0.946979105074,[[280587, 14725], [16570, 278357]],0.98737373989"
$ws4.Range("A11").Value = "2017-08-04 09:29:06.625457 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A12").Value = "1,Test,0.946240603415,[[93263, 4918], [5659, 92907]],0.987243008302"
$ws4.Range("A14").Value = "2017-08-04 09:30:19.248778 Model Training Complete"
$ws4.Range("A15").Value = "2017-08-04 09:30:44.300483 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A16").Value = "2,Train,0.946553944158,[[280029, 14867], [16679, 278665]],0.987345656499"
$ws4.Range("A17").Value = "2017-08-04 09:30:53.231548 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A18").Value = "2,Test,0.947251786567,[[93724, 4873], [5505, 92644]],0.987606348865"
$ws4.Range("A20").Value = "2017-08-04 09:32:05.876631 Model Training Complete"
$ws4.Range("A21").Value = "2017-08-04 09:32:30.917370 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A22").Value = "3,Train,0.946770805096,[[280581, 14421], [16997, 278241]],0.987195909896"
$ws4.Range("A23").Value = "2017-08-04 09:32:39.865750 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("A24").Value = "3,Test,0.94700273449,[[93613, 4878], [5549, 92706]],0.987080938699"
$ws4.Range("J2").Value = "2017-08-04 10:13:34.796964 Model Training Complete"
$ws4.Range("J3").Value = "2017-08-04 10:13:53.703949 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J4").Value = "0,Train,0.961879848671,[[282815, 12454], [10046, 284924]],0.993794558444"
$ws4.Range("J5").Value = "2017-08-04 10:14:00.366457 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J6").Value = "0,Test,0.961447950922,[[94083, 4141], [3444, 95079]],0.993665645289"
$ws4.Range("J7").Value = "2017-08-04 10:14:52.359778 Model Training Complete"
$ws4.Range("J8").Value = "2017-08-04 10:15:12.756062 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J9").Value = "1,Train,0.961842575635,[[282938, 12374], [10148, 284779]],0.993839455224"
$ws4.Range("J10").Value = "2017-08-04 10:15:19.537915 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J11").Value = "1,Test,0.961595348341,[[94050, 4131], [3425, 95141]],0.993599274494"
$ws4.Range("J12").Value = "2017-08-04 10:16:10.275047 Model Training Complete"
$ws4.Range("J13").Value = "2017-08-04 10:16:30.884012 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J14").Value = "2,Train,0.961773177013,[[282402, 12494], [10069, 285275]],0.993760912458"
$ws4.Range("J15").Value = "2017-08-04 10:16:37.557245 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J16").Value = "2,Test,0.961269860632,[[94369, 4228], [3392, 94757]],0.993692358821"
$ws4.Range("J17").Value = "2017-08-04 10:17:28.028609 Model Training Complete"
$ws4.Range("J18").Value = "2017-08-04 10:17:48.449445 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J19").Value = "3,Train,0.961834169152,[[282348, 12654], [9873, 285365]],0.993835736143"
$ws4.Range("J20").Value = "2017-08-04 10:17:55.164115 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J21").Value = "3,Test,0.961036056642,[[94160, 4331], [3335, 94920]],0.993568373451"
$ws4.Range("J1").Value = "2017-08-04 10:12:50.654994 Training XGBoost classifier, [0.6, 5, 1, 0.6]"
$ws4.Range("J25").Value = "2017-08-04 10:59:23.929414 Model Training Complete"
$ws4.Range("J26").Value = "2017-08-04 10:59:42.940074 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J27").Value = "0,Train,0.961879848671,[[282815, 12454], [10046, 284924]],0.993794558444"
$ws4.Range("J28").Value = "2017-08-04 10:59:49.741235 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J29").Value = "0,Test,0.961447950922,[[94083, 4141], [3444, 95079]],0.993665645289"
$ws4.Range("J30").Value = "2017-08-04 11:00:42.941637 Model Training Complete"
$ws4.Range("J31").Value = "2017-08-04 11:01:03.437199 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J32").Value = "1,Train,0.961842575635,[[282938, 12374], [10148, 284779]],0.993839455224"
$ws4.Range("J33").Value = "2017-08-04 11:01:10.150712 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J34").Value = "1,Test,0.961595348341,[[94050, 4131], [3425, 95141]],0.993599274494"
$ws4.Range("J35").Value = "2017-08-04 11:02:02.329446 Model Training Complete"
$ws4.Range("J36").Value = "2017-08-04 11:02:22.911232 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J37").Value = "2,Train,0.961773177013,[[282402, 12494], [10069, 285275]],0.993760912458"
$ws4.Range("J38").Value = "2017-08-04 11:02:29.702509 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J39").Value = "2,Test,0.961269860632,[[94369, 4228], [3392, 94757]],0.993692358821"
$ws4.Range("J40").Value = "2017-08-04 11:03:22.223894 Model Training Complete"
$ws4.Range("J41").Value = "2017-08-04 11:03:42.815590 KF_Index,Train,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J42").Value = "3,Train,0.961834169152,[[282348, 12654], [9873, 285365]],0.993835736143"
$ws4.Range("J43").Value = "2017-08-04 11:03:49.540223 KF_Index,Test,Accuracy,Confusion_Matrix,AUC"
$ws4.Range("J44").Value = "3,Test,0.961036056642,[[94160, 4331], [3335, 94920]],0.993568373451"
$ws4.Range("J24").Value = "2017-08-04 10:58:40.708023 Training XGBoost classifier, [0.6, 5, 1, 0.6]"

# --- View state -------------------------------------------------------
# Sheet3: move the visible window / selection (tabSelected moves off Sheet3
# once Sheet4 is activated below).
$ws3.Activate()
$ws3.Range("W4").Select()

# Sheet4 becomes the active/visible sheet with its own selection.
$ws4.Activate()
$ws4.Range("J25").Select()
